$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row for C9 (AC coupling capacitor) above the old row 4 (D1) ---
# This shifts D1..(old row20) down by one row, matching the target layout.
$ws.Rows.Item(4).Insert()

# --- Row 2: C1, C2 -- 22pF capacitor, now with Manufacturer/MPN filled in ---
$ws.Range("B2").Value = "22pF 0805 C0G Capacitor"
$ws.Range("C2").Value = "Kemet"
$ws.Range("D2").Value = "C0805C220J5GACTU"

# --- Row 3: C3-C7 -> C3-C6, C8 -- 1uF capacitor, now with Manufacturer/MPN filled in ---
$ws.Range("A3").Value = "C3-C6, C8"
$ws.Range("B3").Value = "1uF 0805 X7R Capacitor"
$ws.Range("C3").Value = "Yageo"
$ws.Range("D3").Value = "CC0805KKX7R7BB105"

# --- Row 4 (new): C9 -- 22uF AC coupling capacitor ---
$ws.Range("A4").Value = "C9"
$ws.Range("B4").Value = "22uF 0805 X6S Capacitor"
$ws.Range("C4").Value = "Murata"
$ws.Range("D4").Value = "GRM21BC80G226ME39L"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.14
$ws.Range("G4").Value = "http://www.digikey.com/product-detail/en/murata-electronics-north-america/GRM21BC80G226ME39L/490-6464-1-ND/3845661"

# --- Row 14 (was row 13, R8 jumper): mark as not-stuffed (NOSTUFF) instead of qty 1 ---
$ws.Range("E14").Value = "NOSTUFF"

# --- Sheet view: selection cell moved ---
$ws.Range("B23").Select() | Out-Null
